$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.14035177230835
$ws.Range("B1").Value = 3.394202709197998
$ws.Range("C1").Value = 6.06645679473877
$ws.Range("D1").Value = 1.746901512145996
$ws.Range("E1").Value = 1.054918885231018
